# Relação dos Conceitos - add the "Recuperar Password" use case to the
# CRUD-matrix table on sheet "Folha1".
#
# The new entry follows exactly the same two-row layout (label row +
# thick-bottom spacer row) and merged-cell pattern as the existing
# "Registar" entry (rows 13:14), with the same C / C / (blank) / C/R
# values, so the simplest, most faithful way to create it is to copy
# that block's formatting/values and then relabel it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "Registar" row-pair (13:14) into the new row-pair
# (29:30): this carries over cell styles, borders and merged cells,
# plus the C / C / <blank> / C/R values used for this entry too.
$ws.Range("D13:H14").Copy($ws.Range("D29:H30"))

# Relabel the new entry.
$ws.Range("D29").Value = "Recuperar Password"

# Match the saved selection/view from the edited workbook.
$ws.Range("K27").Select()
